# Crosscheck workbook: add the 2010 crosscheck column (H) and warm up the
# selection for the 2021/2010 data-generation pass.
#
# Commit: "added crosscheck for 2010 and warm up for data generation 2021/2010"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New column H: "2010" crosscheck, mirroring the existing 2021 (F)
#        and 2007 (G) crosscheck columns ------------------------------------

# Copy the formatting of the "2021" header cell (bold, boxed, centered) onto
# the whole new column range so it reuses the same cell style the other two
# year-header/x-mark columns use, instead of Excel minting a brand-new style.
$ws.Range("F2").Copy()
$ws.Range("H2:H30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Year header
$ws.Range("H2").Value = 2010

# Rows 3-5: special "X " (trailing space) mark
$ws.Range("H3").Value = "X "
$ws.Range("H4").Value = "X "
$ws.Range("H5").Value = "X "

# Rows 6-29: plain "X" mark, same as columns F/G
for ($r = 6; $r -le 29; $r++) {
    $ws.Range("H$r").Value = "X"
}

# Row 30 stays blank (style only), matching F30/G30.

# --- 2. Warm up selection for next data-generation pass ---------------------
$ws.Range("M10").Select()
